$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.485.95"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.938.09"
$ws.Range("E3").Value = "  +4.59%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.59"
$ws.Range("E5").Value = "  +3.07%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2881"
$ws.Range("E8").Value = "  +4.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06645"
$ws.Range("E9").Value = "  +4.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.20"
$ws.Range("E10").Value = "  +6.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.19"
$ws.Range("E11").Value = "  +26.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.925.46"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07620"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("E14").Value = "  +3.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6644"
$ws.Range("E15").Value = "  +6.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "308.30"
$ws.Range("E16").Value = "  +25.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.499.41"

$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007593"
$ws.Range("E20").Value = "  +3.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.177.26"
$ws.Range("E21").Value = "  +2.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.304"
$ws.Range("E22").Value = "  +7.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.313"
$ws.Range("E24").Value = "  +6.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.336"
$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.89"
$ws.Range("E26").Value = "  +2.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.29"
$ws.Range("E27").Value = "  +12.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.059"
$ws.Range("E28").Value = "  +9.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1110"
$ws.Range("E29").Value = "  +7.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.112"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.944"
$ws.Range("E32").Value = "  +2.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05031"
$ws.Range("E33").Value = "  +4.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7428"
$ws.Range("E34").Value = "  +6.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("E35").Value = "  +2.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.763"
$ws.Range("E36").Value = "  +2.35%  "

$ws.Range("E37").Value = "  +3.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.055"
$ws.Range("E39").Value = "  +3.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8820"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.81"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("E42").Value = "  +11.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.817"
$ws.Range("E43").Value = "  +5.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4193"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.304"
$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.246"
$ws.Range("E47").Value = "  +7.48%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1216"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.94"
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05625"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3862"
$ws.Range("E51").Value = "  +4.53%  "
